$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F ("dSF") values to re-pulled data per commit message
# "repull data, push all data, mean calculation"
$ws.Range("F3").Value = -3
$ws.Range("F6").Value = -4
$ws.Range("F8").Value = -3
$ws.Range("F9").Value = -3
$ws.Range("F10").Value = 0
$ws.Range("F13").Value = -3
$ws.Range("F15").Value = -10
$ws.Range("F18").Value = -1
$ws.Range("F20").Value = -5
$ws.Range("F21").Value = 1
$ws.Range("F22").Value = -7
$ws.Range("F25").Value = -4
